$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new "Ewan Marsh" user record as row 33
$ws.Range("A33").Value = 110032
$ws.Range("B33").Value = 9317596770
$ws.Range("C33").Value = "Ewan Marsh"
$ws.Range("D33").Value = "ewan.marsh@xyz.com"
$ws.Range("E33").Value = 818876433
$ws.Range("F33").Value = "ACT"
$ws.Range("G33").Value = "eng"
$ws.Range("H33").Value = "PWD"
$ws.Range("I33").Value = $true
$ws.Range("I33").HorizontalAlignment = -4131
$ws.Range("J33").Value = "superadmin"
$ws.Range("K33").Value = "now()"

# Reflect the author's final viewport/selection state: full-column
# selection starting at L1 (no particular scrolled position)
$ws.Range("L1:XFD1048576").Select() | Out-Null

# Bump the print quality (vertical DPI) on the page setup
$ws.PageSetup.VerticalDpi = 300

Write-Output "row 33 added"
